# "Results qualifier 2 RR vs RCB" — fill in match results for Qualifier 2
# (RR vs RCB) on Sheet1 row 85, propagate the Winner-Prediction coin counts
# on rows 96-102, and log the two new matches' rank-prediction tables on
# Sheet2.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# ---------------------------------------------------------------------
# Sheet1 row 85 ("Qualifier 2"): enter each player's points for the match.
# The D/G/J/M/P/S/V columns are VLOOKUP(RANK(...)) formulas already in
# place; they recompute automatically once the raw scores are entered.
# ---------------------------------------------------------------------
$ws1.Range("E85").Value = 50
$ws1.Range("H85").Value = 70
$ws1.Range("K85").Value = 60
$ws1.Range("N85").Value = 40
$ws1.Range("Q85").Value = 100
$ws1.Range("T85").Value = 0
$ws1.Range("W85").Value = 80

# ---------------------------------------------------------------------
# Winner Prediction - Coins table (rows 96-102): column G holds the coins
# earned from the Qualifier 2 winner prediction for each player.
# ---------------------------------------------------------------------
$ws1.Range("G96").Value = 3
$ws1.Range("G97").Value = 0
$ws1.Range("G98").Value = 0
$ws1.Range("G99").Value = 3
$ws1.Range("G100").Value = 11
$ws1.Range("G101").Value = 0
$ws1.Range("G102").Value = 25

# ---------------------------------------------------------------------
# Sheet2: log the two new matches' rank predictions (Qualifier 2 and
# Finals), matching the layout of the existing blocks above them.
# ---------------------------------------------------------------------

# Qualifier 2 RR vs RCB - header + rank columns, copied from the template
# block (rows 6-13) so the styling/row-height matches, then the values
# are overwritten with this match's data.
$ws2.Range("E6:G13").Copy($ws2.Range("E29"))
$ws2.Range("I6:J7").Copy($ws2.Range("I29"))

$ws2.Range("E29").Value = "Qualifier 2 RR vs RCB"
$ws2.Range("F29").Value = "Predictions Rank 1"
$ws2.Range("G29").Value = "Predictions Rank 2"
$ws2.Range("I29").Value = "Rank 1"
$ws2.Range("J29").Value = "Sundar"

$ws2.Range("E30").Value = "Jaya"
$ws2.Range("F30").Value = "Jaya"
$ws2.Range("G30").Value = "Sundar"
$ws2.Range("I30").Value = "Rank 2"
$ws2.Range("J30").Value = "Vicky"

$ws2.Range("E31").Value = "Justin"
$ws2.Range("F31").Value = "Justin"
$ws2.Range("G31").Value = "Sibi"

$ws2.Range("E32").Value = "Ram"
$ws2.Range("F32").Value = "Ram"
$ws2.Range("G32").Value = "Justin"

$ws2.Range("E33").Value = "Sibi"
$ws2.Range("F33").Value = "Sibi"
$ws2.Range("G33").Value = "Sundar"

$ws2.Range("E34").Value = "Sundar"
$ws2.Range("F34").Value = "Sundar"
$ws2.Range("G34").Value = "Sibi"

$ws2.Range("E35").Value = "Upili"
$ws2.Range("F35").Value = "Upili"
$ws2.Range("G35").Value = "Justin"

$ws2.Range("E36").Value = "Vikcy"
$ws2.Range("F36").Value = "Sundar"
$ws2.Range("G36").Value = "Vicky"

# Finals GT vs RR
$ws2.Range("E6:G13").Copy($ws2.Range("E40"))
$ws2.Range("I6:I7").Copy($ws2.Range("I40"))

$ws2.Range("E40").Value = "Finals GT vs RR"
$ws2.Range("F40").Value = "Predictions Rank 1"
$ws2.Range("G40").Value = "Predictions Rank 2"
$ws2.Range("I40").Value = "Rank 1"
$ws2.Range("J40").ClearContents()

$ws2.Range("E41").Value = "Jaya"
$ws2.Range("F41").Value = "Jaya"
$ws2.Range("G41").Value = "Sundar"
$ws2.Range("I41").Value = "Rank 2"
$ws2.Range("J41").ClearContents()

$ws2.Range("E42").Value = "Justin"
$ws2.Range("F42").Value = "Justin"
$ws2.Range("G42").Value = "Sundar"

$ws2.Range("E43").Value = "Ram"
$ws2.Range("F43").Value = "Sundar"
$ws2.Range("G43").Value = "Sibi"

$ws2.Range("E44").Value = "Sibi"
$ws2.Range("F44").Value = "Sibi"
$ws2.Range("G44").Value = "Sundar"

$ws2.Range("E45").Value = "Sundar"
$ws2.Range("F45").Value = "Sundar"
$ws2.Range("G45").Value = "Vicky"

$ws2.Range("E46").Value = "Upili"
$ws2.Range("F46").Value = "Upili"
$ws2.Range("G46").Value = "Upili"

$ws2.Range("E47").Value = "Vikcy"
$ws2.Range("F47").Value = "Ram"
$ws2.Range("G47").Value = "Vicky"

# ---------------------------------------------------------------------
# View bits: column U got a touch narrower, Sheet2 becomes the active tab.
# ---------------------------------------------------------------------
$ws1.Columns.Item(21).ColumnWidth = 12.25

$ws2.Activate()
$ws2.Range("J41").Select()
